$wb = $excel.ActiveWorkbook

# --- EventsTestData: fix "Due date" -> "Due Date" (new shared string created first) ---
$wsEvents = $wb.Worksheets.Item("EventsTestData")
$wsEvents.Range("C5").Value = "Due Date"

# --- OrganizationsTestData: TC_BWF_02 row (row 5) Status Pass -> Fail (new shared string created second) ---
$wsOrg = $wb.Worksheets.Item("OrganizationsTestData")
$wsOrg.Range("E5").Value = "Fail"

# --- LeadsTestData: TC_BWF_02 row (row 6) Status blank -> Pass ---
$wsLeads = $wb.Worksheets.Item("LeadsTestData")
$wsLeads.Range("E6").Value = "Pass"
$wsLeads.Range("E6").Style = "Normal"

# --- EventsTestData: fill Status on row 2, move selection ---
$wsEvents.Range("E2").Value = "Pass"
$wsEvents.Range("E2").Style = "Normal"
$wsEvents.Activate()
$wsEvents.Range("C5").Select()
